# Automatic update of files.
#
# 1. Column C ("Förändrad") on every data row (rows 2-32) moves from
#    45208 to 45212.
# 2. For the first four data rows (2-5) the link-building formulas in
#    columns S, T, V, W, X, Y get a descriptive suffix appended to the
#    filename (and, for column Y, the target folder name itself changed
#    to "ti,llsynsmail").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. "Förändrad" date column -------------------------------------------------
$ws.Range("C2:C32").Value = 45212

# --- 2. Hyperlink formulas for rows 2-5 -----------------------------------------
$caseIds = @{
    2 = "A 30234-2023"
    3 = "A 33548-2023"
    4 = "A 33550-2023"
    5 = "A 30241-2023"
}

# column -> folder, suffix appended before extension, extension
$linkCols = @(
    @{ Col = "S"; Folder = "artfynd";       Suffix = " artfynd";              Ext = ".xlsx" },
    @{ Col = "T"; Folder = "kartor";        Suffix = " karta";                Ext = ".png"  },
    @{ Col = "V"; Folder = "klagomål";      Suffix = " fsc-klagomål";         Ext = ".docx" },
    @{ Col = "W"; Folder = "klagomålsmail"; Suffix = " fsc-klagomål mail";    Ext = ".docx" },
    @{ Col = "X"; Folder = "tillsyn";       Suffix = " tillsynsbegäran";      Ext = ".docx" },
    @{ Col = "Y"; Folder = "ti,llsynsmail"; Suffix = " tillsynsbegäran mail"; Ext = ".docx" }
)

foreach ($row in $caseIds.Keys) {
    $caseId = $caseIds[$row]
    foreach ($link in $linkCols) {
        $url = "https://klasma.github.io/LoggingDetectiveFiles/Logging_2039/" + $link.Folder + "/" + $caseId + $link.Suffix + $link.Ext
        $formula = '=HYPERLINK("' + $url + '", "' + $caseId + '")'
        $ws.Range($link.Col + $row).Formula = $formula
    }
}
